$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header narrative text (Volume/Number and report date range)
$ws.Range("A8").Value = "Volume 32   Number  42"
$ws.Range("C9").Value = "Report Covering the Week  10/13/2025  Through  10/19/2025"

# --- Cells that change type (number <-> dash/***.* text) ---
# Copy value+format from a same-shaped template cell so the shared-string /
# number style (e.g. style 13 text vs style 14/15 numeric) matches exactly.
$ws.Range("C23").Copy()
$ws.Range("C16").PasteSpecial(-4163)
$ws.Range("C23").Copy()
$ws.Range("C16").PasteSpecial(-4122)

$ws.Range("C19").Copy()
$ws.Range("C17").PasteSpecial(-4163)
$ws.Range("C19").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = 1

$ws.Range("C23").Copy()
$ws.Range("C18").PasteSpecial(-4163)
$ws.Range("C23").Copy()
$ws.Range("C18").PasteSpecial(-4122)

$ws.Range("C23").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("C23").Copy()
$ws.Range("D18").PasteSpecial(-4122)

$ws.Range("E23").Copy()
$ws.Range("E18").PasteSpecial(-4163)
$ws.Range("E23").Copy()
$ws.Range("E18").PasteSpecial(-4122)

$ws.Range("C23").Copy()
$ws.Range("C22").PasteSpecial(-4163)
$ws.Range("C23").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("C19").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("C19").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = 2

$ws.Range("E19").Copy()
$ws.Range("E28").PasteSpecial(-4163)
$ws.Range("E19").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = -50

$ws.Range("C19").Copy()
$ws.Range("G28").PasteSpecial(-4163)
$ws.Range("C19").Copy()
$ws.Range("G28").PasteSpecial(-4122)
$ws.Range("G28").Value = 2

$ws.Range("E19").Copy()
$ws.Range("H28").PasteSpecial(-4163)
$ws.Range("E19").Copy()
$ws.Range("H28").PasteSpecial(-4122)
$ws.Range("H28").Value = 200

$ws.Range("C23").Copy()
$ws.Range("G33").PasteSpecial(-4163)
$ws.Range("C23").Copy()
$ws.Range("G33").PasteSpecial(-4122)

$ws.Range("E23").Copy()
$ws.Range("H33").PasteSpecial(-4163)
$ws.Range("E23").Copy()
$ws.Range("H33").PasteSpecial(-4122)

# --- Plain value updates (style/type unchanged) ---
$ws.Range("D15").Value = 2
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = -80
$ws.Range("J15").Value = 12
$ws.Range("K15").Value = 0
$ws.Range("N15").Value = -52
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = -100
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = -57.142857142857
$ws.Range("J16").Value = 57
$ws.Range("K16").Value = -14.035087719298
$ws.Range("L16").Value = -14.035087719298
$ws.Range("M16").Value = -40.963855421686
$ws.Range("N16").Value = -89.207048458149
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 6
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 111
$ws.Range("J17").Value = 97
$ws.Range("K17").Value = 14.432989690721
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 32.142857142857
$ws.Range("N17").Value = -51.948051948051
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 40
$ws.Range("M18").Value = -69.651741293532
$ws.Range("N18").Value = -93.153759820426
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 29
$ws.Range("H19").Value = -13.793103448275
$ws.Range("I19").Value = 302
$ws.Range("J19").Value = 304
$ws.Range("K19").Value = -0.657894736842
$ws.Range("L19").Value = -14.447592067988
$ws.Range("M19").Value = 18.431372549019
$ws.Range("N19").Value = -18.378378378378
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -80
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = -37.5
$ws.Range("I20").Value = 107
$ws.Range("J20").Value = 162
$ws.Range("K20").Value = -33.95061728395
$ws.Range("L20").Value = 2.884615384615
$ws.Range("M20").Value = -17.054263565891
$ws.Range("N20").Value = -93.069948186528
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = -50
$ws.Range("F21").Value = 52
$ws.Range("G21").Value = 71
$ws.Range("H21").Value = -26.760563380281
$ws.Range("I21").Value = 642
$ws.Range("J21").Value = 708
$ws.Range("K21").Value = -9.322033898305
$ws.Range("L21").Value = -9.957924263674
$ws.Range("M21").Value = -15.415019762845
$ws.Range("N21").Value = -81.756180733162
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 9
$ws.Range("K22").Value = -33.333333333333
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = -23.809523809523
$ws.Range("F24").Value = 57
$ws.Range("G24").Value = 93
$ws.Range("H24").Value = -38.709677419354
$ws.Range("I24").Value = 791
$ws.Range("J24").Value = 1148
$ws.Range("K24").Value = -31.097560975609
$ws.Range("L24").Value = -36.209677419354
$ws.Range("M24").Value = -18.453608247422
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = -81.25
$ws.Range("F25").Value = 23
$ws.Range("H25").Value = -58.928571428571
$ws.Range("I25").Value = 364
$ws.Range("J25").Value = 715
$ws.Range("K25").Value = -49.090909090909
$ws.Range("L25").Value = -44.680851063829
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 12.5
$ws.Range("F26").Value = 25
$ws.Range("G26").Value = 27
$ws.Range("H26").Value = -7.407407407407
$ws.Range("I26").Value = 317
$ws.Range("J26").Value = 288
$ws.Range("K26").Value = 10.069444444444
$ws.Range("L26").Value = 6.020066889632
$ws.Range("M26").Value = 1.92926045016
$ws.Range("D27").Value = 2
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -33.333333333333
$ws.Range("J27").Value = 16
$ws.Range("K27").Value = 18.75
$ws.Range("I28").Value = 37
$ws.Range("J28").Value = 44
$ws.Range("K28").Value = -15.90909090909
$ws.Range("L28").Value = 60.869565217391
$ws.Range("L31").Value = -42.857142857142
